$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for new columns I (I0) and J (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (style) from an existing header cell (H1) onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill column J (IF) with the same values as column H (IP) for data rows 2-30
for ($r = 2; $r -le 30; $r++) {
    $hval = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 10).Value = $hval
}

# Fill column I (I0) with 1 for all data rows 2-30
$ws.Range("I2:I30").Value = 1

# Row 3 is a special case in the source data: I3=5, J3=7 (not following the
# general I=1 / J=H pattern used by every other row)
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 7
